$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 24 cell updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2221.2856
$ws.Range("I28").Value = 2209.8
$ws.Range("K28").Value = 2209.8
$ws.Range("M28").Value = -1724.8
$ws.Range("H33").Value = 49.96154
$ws.Range("I33").Value = 49.96154
$ws.Range("K33").Value = 49.96154
$ws.Range("M33").Value = 179.03846
$ws.Range("H62").Value = 799
$ws.Range("I62").Value = 799
$ws.Range("K62").Value = 799
$ws.Range("M62").Value = -175
$ws.Range("H65").Value = 799
$ws.Range("I65").Value = 799
$ws.Range("K65").Value = 3995
$ws.Range("M65").Value = -875
$ws.Range("H92").Value = 317.63635
$ws.Range("I92").Value = 316
$ws.Range("K92").Value = 316
$ws.Range("M92").Value = 932
$ws.Range("H100").Value = 2721.75
$ws.Range("I100").Value = 1943.5
$ws.Range("K100").Value = 1943.5
$ws.Range("M100").Value = -1402.5

# --- Sheet ARM: 16 cell updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 20000
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H20").Value = 20000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet BSM: 8 cell updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 422.8
$ws.Range("J12").Value = 500
$ws.Range("L12").Value = 500
$ws.Range("N12").Value = -836
$ws.Range("H107").Value = 1758.6666
$ws.Range("I107").Value = 1758.6666
$ws.Range("K107").Value = 1758.6666
$ws.Range("M107").Value = 161.3334

# --- Sheet CRP: 38 cell updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2981.4
$ws.Range("I5").Value = 3226.75
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 3226.75
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -3114.75
$ws.Range("N5").Value = -2224
$ws.Range("H74").Value = 98799.664
$ws.Range("J74").Value = 98799.664
$ws.Range("L74").Value = 98799.664
$ws.Range("N74").Value = -100547.664
$ws.Range("H77").Value = 98799.664
$ws.Range("J77").Value = 98799.664
$ws.Range("L77").Value = 296398.992
$ws.Range("N77").Value = -305134.992
$ws.Range("H86").Value = 8105.6665
$ws.Range("J86").Value = 6949
$ws.Range("L86").Value = 6949
$ws.Range("N86").Value = -9195
$ws.Range("H89").Value = 8105.6665
$ws.Range("J89").Value = 6949
$ws.Range("L89").Value = 34745
$ws.Range("N89").Value = -45977
$ws.Range("H95").Value = 32924.6
$ws.Range("J95").Value = 32924.6
$ws.Range("L95").Value = 32924.6
$ws.Range("N95").Value = -38416.6
$ws.Range("H107").Value = 180.88889
$ws.Range("I107").Value = 208.38461
$ws.Range("K107").Value = 208.38461
$ws.Range("M107").Value = 1711.61539
$ws.Range("H122").Value = 1257.1
$ws.Range("I122").Value = 1393.8572
$ws.Range("J122").Value = 938
$ws.Range("K122").Value = 4181.571599999999
$ws.Range("L122").Value = 2814
$ws.Range("M122").Value = -1731.571599999999
$ws.Range("N122").Value = -7714

# --- Sheet CUL: 46 cell updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 40843.92
$ws.Range("I4").Value = 1061.8182
$ws.Range("J4").Value = 113777.78
$ws.Range("K4").Value = 3185.4546
$ws.Range("L4").Value = 341333.34
$ws.Range("M4").Value = -3073.4546
$ws.Range("N4").Value = -341557.34
$ws.Range("H54").Value = 6000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 6000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 18000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -19118
$ws.Range("H61").Value = 432.85715
$ws.Range("I61").Value = 88.333336
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 265.000008
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -50.00000799999998
$ws.Range("N61").Value = -7930
$ws.Range("H68").Value = 576
$ws.Range("I68").Value = 590
$ws.Range("J68").Value = 566.6667
$ws.Range("K68").Value = 1770
$ws.Range("L68").Value = 1700.0001
$ws.Range("M68").Value = -959
$ws.Range("N68").Value = -3322.0001
$ws.Range("H71").Value = 576
$ws.Range("I71").Value = 590
$ws.Range("J71").Value = 566.6667
$ws.Range("K71").Value = 5310
$ws.Range("L71").Value = 5100.0003
$ws.Range("M71").Value = -1254
$ws.Range("N71").Value = -13212.0003
$ws.Range("H109").Value = 975.125
$ws.Range("I109").Value = 1501
$ws.Range("K109").Value = 4503
$ws.Range("M109").Value = -3463
$ws.Range("H140").Value = 965
$ws.Range("I140").Value = 430
$ws.Range("J140").Value = 1500
$ws.Range("K140").Value = 1290
$ws.Range("L140").Value = 4500
$ws.Range("M140").Value = 3890
$ws.Range("N140").Value = -14860

# --- Sheet GSM: 19 cell updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 403.5
$ws.Range("I9").Value = 403.5
$ws.Range("K9").Value = 403.5
$ws.Range("M9").Value = -233.5
$ws.Range("H13").Value = 128.33333
$ws.Range("J13").Value = 150
$ws.Range("L13").Value = 150
$ws.Range("N13").Value = -428
$ws.Range("H17").Value = 95
$ws.Range("I17").Value = 95
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 95
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 73
$ws.Range("N17").ClearContents()
$ws.Range("H52").Value = 20000
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# --- Sheet LTW: 27 cell updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 301
$ws.Range("I9").Value = 325
$ws.Range("K9").Value = 325
$ws.Range("M9").Value = -101
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H31").Value = 3605
$ws.Range("I31").Value = 995
$ws.Range("J31").Value = 4040
$ws.Range("K31").Value = 995
$ws.Range("L31").Value = 4040
$ws.Range("M31").Value = -747
$ws.Range("N31").Value = -4536
$ws.Range("H58").Value = 33773.25
$ws.Range("I58").Value = 21546.5
$ws.Range("K58").Value = 21546.5
$ws.Range("M58").Value = -21286.5
$ws.Range("H61").Value = 3230.625
$ws.Range("I61").Value = 3535
$ws.Range("K61").Value = 3535
$ws.Range("M61").Value = -3333
$ws.Range("H113").Value = 3230.625
$ws.Range("I113").Value = 3535
$ws.Range("K113").Value = 3535
$ws.Range("M113").Value = -1365

# --- Sheet WVR: 21 cell updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2449.3333
$ws.Range("I81").Value = 2449.3333
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4898.6666
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -3837.6666
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 2449.3333
$ws.Range("I84").Value = 2449.3333
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 24493.333
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -19189.333
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 4268.357
$ws.Range("I113").Value = 365.77777
$ws.Range("J113").Value = 11293
$ws.Range("K113").Value = 1097.33331
$ws.Range("L113").Value = 33879
$ws.Range("M113").Value = 1072.66669
$ws.Range("N113").Value = -38219
